# Fruta / hortaliza, semanal
# Insert a new weekly observation row at row 117 (pushing the existing
# rows 117-142 down to 118-143) and populate the new row with the
# latest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 117, shifting rows 117:142 -> 118:143
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new weekly record
$ws.Range("A117").Value = 10
$ws.Range("B117").Value = "Vega Modelo de Temuco"
$ws.Range("C117").Value = "La Araucanía"
$ws.Range("D117").Value = 44754
$ws.Range("E117").Value = 9
$ws.Range("F117").Value = 100112031
$ws.Range("G117").Value = "Poroto verde"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 5
$ws.Range("K117").Value = 40000
$ws.Range("L117").Value = 40000
$ws.Range("M117").Value = 40000
$ws.Range("N117").Value = "$/malla 25 kilos"
$ws.Range("O117").Value = "Provincia de Limarí"
$ws.Range("P117").Value = 1600
$ws.Range("Q117").Value = 25
$ws.Range("R117").Value = "Hortaliza"
